# Worked on temporal resolution
# Extend the "Demand" sheet's single time-step row (row 3) into a full
# 12-step time series (rows 3-14), matching the change in demand data
# for commodity Co in site Sit across the additional time steps.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demand")

# Make "Demand" the active sheet (it becomes the tab shown/selected).
$ws.Activate()

# Copy the formatting of the existing data cell (B3, plain/default style)
# down across the new rows so the newly written values keep the same
# (unstyled) look as the original entry - mirrors a fill-down in Excel.
$ws.Range("B3").Copy()
$ws.Range("B4:B14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column A: time step index 1..12
# Column B: demand value (MWh) for EU27.Elec - same figure repeated
# across every time step (451750000), replacing the old single lump
# value of 5421000000 in B3 (=12 * 451750000).
for ($row = 3; $row -le 14; $row++) {
    $ws.Cells.Item($row, 1).Value = $row - 2
    $ws.Cells.Item($row, 2).Value = 451750000
}

# Leave the new block selected, like it was right after entering it.
$ws.Range("B3:B14").Select()
